$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("D12").Value = "keras-nightly 패키지"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/05/18/keras-nightly-%ed%8c%a8%ed%82%a4%ec%a7%80/"

# Row 28
$ws.Range("D28").Value = "[Null space control] Null space control이란?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/126"

# Row 29
$ws.Range("D29").Value = "[만화] 인턴일기 28~33"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-5/"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Text-to-SQL"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1482&mod=document&pageid=1"

# Row 51
$ws.Range("D51").Value = "[MariaDB] SELECT 한 것을 INSERT 하기, INSERT INTO ... SELECT문"
$ws.Range("E51").Value = "https://bskyvision.com/1188"
